$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 7576918   # H17: was 3247928.2
$ws.Cells.Item(17, 9).Value = 0   # I17: was 1139
$ws.Cells.Item(17, 10).Value = 7576918   # J17: was 4133416.2
$ws.Cells.Item(17, 11).Value = 0   # K17: was 3417
$ws.Cells.Item(17, 12).Value = 22730754   # L17: was 12400248.6
$ws.Cells.Item(17, 13).Value = ""   # M17: was -3249, clear cell
$ws.Cells.Item(17, 14).Value = -22731090   # N17: was -12400584.6
$ws.Cells.Item(46, 8).Value = 203599   # H46: was 1000000
$ws.Cells.Item(46, 9).Value = 3000   # I46: was 0
$ws.Cells.Item(46, 10).Value = 253748.75   # J46: was 1000000
$ws.Cells.Item(46, 11).Value = 9000   # K46: was 0
$ws.Cells.Item(46, 12).Value = 761246.25   # L46: was 3000000
$ws.Cells.Item(46, 13).Value = -8881   # M46: was None
$ws.Cells.Item(46, 14).Value = -761484.25   # N46: was -3000238
$ws.Cells.Item(60, 8).Value = 203599   # H60: was 1000000
$ws.Cells.Item(60, 9).Value = 3000   # I60: was 0
$ws.Cells.Item(60, 10).Value = 253748.75   # J60: was 1000000
$ws.Cells.Item(60, 11).Value = 9000   # K60: was 0
$ws.Cells.Item(60, 12).Value = 761246.25   # L60: was 3000000
$ws.Cells.Item(60, 13).Value = -8516   # M60: was None
$ws.Cells.Item(60, 14).Value = -762214.25   # N60: was -3000968
$ws.Cells.Item(62, 8).Value = 4927.5386   # H62: was 4934.5713
$ws.Cells.Item(62, 9).Value = 4921.5   # I62: was 4929.5386
$ws.Cells.Item(62, 11).Value = 4921.5   # K62: was 4929.5386
$ws.Cells.Item(62, 13).Value = -4297.5   # M62: was -4305.5386
$ws.Cells.Item(65, 8).Value = 4927.5386   # H65: was 4934.5713
$ws.Cells.Item(65, 9).Value = 4921.5   # I65: was 4929.5386
$ws.Cells.Item(65, 11).Value = 24607.5   # K65: was 24647.693
$ws.Cells.Item(65, 13).Value = -21487.5   # M65: was -21527.693
$ws.Cells.Item(70, 8).Value = 40001348   # H70: was 20001174
$ws.Cells.Item(70, 9).Value = 200000000   # I70: was 40000800
$ws.Cells.Item(70, 10).Value = 1687.25   # J70: was 1549.8
$ws.Cells.Item(70, 11).Value = 600000000   # K70: was 120002400
$ws.Cells.Item(70, 12).Value = 5061.75   # L70: was 4649.4
$ws.Cells.Item(70, 13).Value = -599999730   # M70: was -120002130
$ws.Cells.Item(70, 14).Value = -5601.75   # N70: was -5189.4
$ws.Cells.Item(73, 8).Value = 40001348   # H73: was 20001174
$ws.Cells.Item(73, 9).Value = 200000000   # I73: was 40000800
$ws.Cells.Item(73, 10).Value = 1687.25   # J73: was 1549.8
$ws.Cells.Item(73, 11).Value = 600000000   # K73: was 120002400
$ws.Cells.Item(73, 12).Value = 5061.75   # L73: was 4649.4
$ws.Cells.Item(73, 13).Value = -599999064   # M73: was -120001464
$ws.Cells.Item(73, 14).Value = -6933.75   # N73: was -6521.4
$ws.Cells.Item(113, 8).Value = 16354.556   # H113: was 13098.833
$ws.Cells.Item(113, 9).Value = 29999   # I113: was 19199.4
$ws.Cells.Item(113, 10).Value = 9532.333000000001   # J113: was 8741.286
$ws.Cells.Item(113, 11).Value = 29999   # K113: was 19199.4
$ws.Cells.Item(113, 12).Value = 9532.333000000001   # L113: was 8741.286
$ws.Cells.Item(113, 13).Value = -26745   # M113: was -15945.4
$ws.Cells.Item(113, 14).Value = -16040.333   # N113: was -15249.286
$ws.Cells.Item(133, 8).Value = 83850.336   # H133: was 84309.60000000001
$ws.Cells.Item(133, 10).Value = 83850.336   # J133: was 84309.60000000001
$ws.Cells.Item(133, 12).Value = 83850.336   # L133: was 84309.60000000001
$ws.Cells.Item(133, 14).Value = -93970.336   # N133: was -94429.60000000001

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 252250   # H2: was 202780
$ws.Cells.Item(2, 9).Value = 1000   # I2: was 2950
$ws.Cells.Item(2, 11).Value = 1000   # K2: was 2950
$ws.Cells.Item(2, 13).Value = -887   # M2: was -2837
$ws.Cells.Item(32, 8).Value = 6444.731   # H32: was 2887.322
$ws.Cells.Item(32, 9).Value = 6302.52   # I32: was 2764.6897
$ws.Cells.Item(32, 11).Value = 6302.52   # K32: was 2764.6897
$ws.Cells.Item(32, 13).Value = -6015.52   # M32: was -2477.6897
$ws.Cells.Item(61, 8).Value = 5343.3335   # H61: was 5467.3477
$ws.Cells.Item(61, 9).Value = 5353.913   # I61: was 5484.0454
$ws.Cells.Item(61, 11).Value = 5353.913   # K61: was 5484.0454
$ws.Cells.Item(61, 13).Value = -5141.913   # M61: was -5272.0454
$ws.Cells.Item(74, 8).Value = 5036.3687   # H74: was 5504.4375
$ws.Cells.Item(74, 9).Value = 1835.0714   # I74: was 1922.5
$ws.Cells.Item(74, 10).Value = 14000   # J74: was 16250.25
$ws.Cells.Item(74, 11).Value = 1835.0714   # K74: was 1922.5
$ws.Cells.Item(74, 12).Value = 14000   # L74: was 16250.25
$ws.Cells.Item(74, 13).Value = -961.0714   # M74: was -1048.5
$ws.Cells.Item(74, 14).Value = -15748   # N74: was -17998.25
$ws.Cells.Item(77, 8).Value = 5036.3687   # H77: was 5504.4375
$ws.Cells.Item(77, 9).Value = 1835.0714   # I77: was 1922.5
$ws.Cells.Item(77, 10).Value = 14000   # J77: was 16250.25
$ws.Cells.Item(77, 11).Value = 9175.357   # K77: was 9612.5
$ws.Cells.Item(77, 12).Value = 70000   # L77: was 81251.25
$ws.Cells.Item(77, 13).Value = -4807.357   # M77: was -5244.5
$ws.Cells.Item(77, 14).Value = -78736   # N77: was -89987.25
$ws.Cells.Item(109, 8).Value = 0   # H109: was 20000
$ws.Cells.Item(109, 10).Value = 0   # J109: was 20000
$ws.Cells.Item(109, 12).Value = 0   # L109: was 20000
$ws.Cells.Item(109, 14).Value = ""   # N109: was -22774, clear cell
$ws.Cells.Item(116, 8).Value = 252250   # H116: was 202780
$ws.Cells.Item(116, 9).Value = 1000   # I116: was 2950
$ws.Cells.Item(116, 11).Value = 1000   # K116: was 2950
$ws.Cells.Item(116, 13).Value = 1294   # M116: was -656
$ws.Cells.Item(122, 8).Value = 487851.94   # H122: was 472278.03
$ws.Cells.Item(122, 9).Value = 4034.9   # I122: was 4768.8887
$ws.Cells.Item(122, 10).Value = 1563000.9   # J122: was 1173541.8
$ws.Cells.Item(122, 11).Value = 12104.7   # K122: was 14306.6661
$ws.Cells.Item(122, 12).Value = 4689002.699999999   # L122: was 3520625.4
$ws.Cells.Item(122, 13).Value = -9654.700000000001   # M122: was -11856.6661
$ws.Cells.Item(122, 14).Value = -4693902.699999999   # N122: was -3525525.4
$ws.Cells.Item(136, 8).Value = 5343.3335   # H136: was 5467.3477
$ws.Cells.Item(136, 9).Value = 5353.913   # I136: was 5484.0454
$ws.Cells.Item(136, 11).Value = 16061.739   # K136: was 16452.1362
$ws.Cells.Item(136, 13).Value = -13511.739   # M136: was -13902.1362

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 252250   # H3: was 202780
$ws.Cells.Item(3, 9).Value = 1000   # I3: was 2950
$ws.Cells.Item(3, 11).Value = 1000   # K3: was 2950
$ws.Cells.Item(3, 13).Value = -886   # M3: was -2836
$ws.Cells.Item(20, 8).Value = 3362.625   # H20: was 2931.1052
$ws.Cells.Item(20, 9).Value = 2487.6   # I20: was 2188.75
$ws.Cells.Item(20, 10).Value = 4821   # J20: was 4203.7144
$ws.Cells.Item(20, 11).Value = 2487.6   # K20: was 2188.75
$ws.Cells.Item(20, 12).Value = 4821   # L20: was 4203.7144
$ws.Cells.Item(20, 13).Value = -2240.6   # M20: was -1941.75
$ws.Cells.Item(20, 14).Value = -5315   # N20: was -4697.7144
$ws.Cells.Item(116, 8).Value = 69000   # H116: was 68475
$ws.Cells.Item(116, 10).Value = 69000   # J116: was 68475
$ws.Cells.Item(116, 12).Value = 69000   # L116: was 68475
$ws.Cells.Item(116, 14).Value = -78178   # N116: was -77653
$ws.Cells.Item(124, 8).Value = 63000   # H124: was 64000
$ws.Cells.Item(124, 10).Value = 63000   # J124: was 64000
$ws.Cells.Item(124, 12).Value = 63000   # L124: was 64000
$ws.Cells.Item(124, 14).Value = -72820   # N124: was -73820

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1208.5454   # H16: was 1228.4
$ws.Cells.Item(16, 9).Value = 1208.5454   # I16: was 1228.4
$ws.Cells.Item(16, 11).Value = 1208.5454   # K16: was 1228.4
$ws.Cells.Item(16, 13).Value = -921.5454   # M16: was -941.4000000000001
$ws.Cells.Item(105, 8).Value = 213679.9   # H105: was 266737.5
$ws.Cells.Item(105, 9).Value = 265949.88   # I105: was 354116.66
$ws.Cells.Item(105, 11).Value = 265949.88   # K105: was 354116.66
$ws.Cells.Item(105, 13).Value = -264202.88   # M105: was -352369.66
$ws.Cells.Item(107, 8).Value = 8481.448   # H107: was 8201.267
$ws.Cells.Item(107, 9).Value = 12014.35   # I107: was 11445.857
$ws.Cells.Item(107, 11).Value = 12014.35   # K107: was 11445.857
$ws.Cells.Item(107, 13).Value = -10094.35   # M107: was -9525.857
$ws.Cells.Item(113, 8).Value = 1208.5454   # H113: was 1228.4
$ws.Cells.Item(113, 9).Value = 1208.5454   # I113: was 1228.4
$ws.Cells.Item(113, 11).Value = 1208.5454   # K113: was 1228.4
$ws.Cells.Item(113, 13).Value = 961.4546   # M113: was 941.5999999999999
$ws.Cells.Item(122, 8).Value = 12415.909   # H122: was 11418.75
$ws.Cells.Item(122, 9).Value = 40666.668   # I122: was 30612.5
$ws.Cells.Item(122, 11).Value = 122000.004   # K122: was 91837.5
$ws.Cells.Item(122, 13).Value = -119550.004   # M122: was -89387.5
$ws.Cells.Item(132, 8).Value = 38706.184   # H132: was 42276.8
$ws.Cells.Item(132, 9).Value = 2529.7778   # I132: was 2471
$ws.Cells.Item(132, 11).Value = 7589.3334   # K132: was 7413
$ws.Cells.Item(132, 13).Value = -5059.3334   # M132: was -4883
$ws.Cells.Item(134, 8).Value = 3563.8667   # H134: was 3321.0588
$ws.Cells.Item(134, 9).Value = 3788.4167   # I134: was 3461.5
$ws.Cells.Item(134, 11).Value = 11365.2501   # K134: was 10384.5
$ws.Cells.Item(134, 13).Value = -8830.250100000001   # M134: was -7849.5

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(58, 8).Value = 3466.5833   # H58: was 3500
$ws.Cells.Item(58, 9).Value = 3299.5   # I58: was 0
$ws.Cells.Item(58, 11).Value = 9898.5   # K58: was 0
$ws.Cells.Item(58, 13).Value = -9770.5   # M58: was None
$ws.Cells.Item(63, 8).Value = 2330.2856   # H63: was 2402
$ws.Cells.Item(63, 9).Value = 1900   # I63: was 0
$ws.Cells.Item(63, 11).Value = 5700   # K63: was 0
$ws.Cells.Item(63, 13).Value = -4951   # M63: was None
$ws.Cells.Item(66, 8).Value = 2330.2856   # H66: was 2402
$ws.Cells.Item(66, 9).Value = 1900   # I66: was 0
$ws.Cells.Item(66, 11).Value = 17100   # K66: was 0
$ws.Cells.Item(66, 13).Value = -13356   # M66: was None
$ws.Cells.Item(114, 8).Value = 5499.6   # H114: was 6166
$ws.Cells.Item(114, 10).Value = 6249.75   # J114: was 7999.5
$ws.Cells.Item(114, 12).Value = 18749.25   # L114: was 23998.5
$ws.Cells.Item(114, 14).Value = -25257.25   # N114: was -30506.5
$ws.Cells.Item(127, 8).Value = 1001   # H127: was 100
$ws.Cells.Item(127, 10).Value = 1001   # J127: was 100
$ws.Cells.Item(127, 12).Value = 3003   # L127: was 300
$ws.Cells.Item(127, 14).Value = -12923   # N127: was -10220
$ws.Cells.Item(131, 8).Value = 15386058   # H131: was 14707346
$ws.Cells.Item(131, 10).Value = 1533.6364   # J131: was 1553.9656
$ws.Cells.Item(131, 12).Value = 4600.9092   # L131: was 4661.8968
$ws.Cells.Item(131, 14).Value = -14680.9092   # N131: was -14741.8968
$ws.Cells.Item(139, 8).Value = 2223825.5   # H139: was 2223947.5
$ws.Cells.Item(139, 9).Value = 2858061.2   # I139: was 2858218.2
$ws.Cells.Item(139, 11).Value = 8574183.600000001   # K139: was 8574654.600000001
$ws.Cells.Item(139, 13).Value = -8569043.600000001   # M139: was -8569514.600000001

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 1652.6666   # H132: was 1699.4286
$ws.Cells.Item(132, 9).Value = 1574.1666   # I132: was 1626.5454
$ws.Cells.Item(132, 11).Value = 4722.4998   # K132: was 4879.6362
$ws.Cells.Item(132, 13).Value = -2192.4998   # M132: was -2349.6362

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 21015.215   # H7: was 20425.068
$ws.Cells.Item(7, 9).Value = 44874.2   # I7: was 38045.082
$ws.Cells.Item(7, 10).Value = 7760.222   # J7: was 7987.4116
$ws.Cells.Item(7, 11).Value = 44874.2   # K7: was 38045.082
$ws.Cells.Item(7, 12).Value = 7760.222   # L7: was 7987.4116
$ws.Cells.Item(7, 13).Value = -44762.2   # M7: was -37933.082
$ws.Cells.Item(7, 14).Value = -7984.222   # N7: was -8211.411599999999
$ws.Cells.Item(16, 8).Value = 6364.143   # H16: was 6591.8335
$ws.Cells.Item(16, 10).Value = 4999   # J16: was 4999.5
$ws.Cells.Item(16, 12).Value = 4999   # L16: was 4999.5
$ws.Cells.Item(16, 14).Value = -5339   # N16: was -5339.5
$ws.Cells.Item(19, 8).Value = 2325   # H19: was 1766.6666
$ws.Cells.Item(46, 8).Value = 4015.2856   # H46: was 4495.5835
$ws.Cells.Item(46, 9).Value = 1125   # I46: was 983.3333
$ws.Cells.Item(46, 10).Value = 4497   # J46: was 5666.3335
$ws.Cells.Item(46, 11).Value = 1125   # K46: was 983.3333
$ws.Cells.Item(46, 12).Value = 4497   # L46: was 5666.3335
$ws.Cells.Item(46, 13).Value = -937   # M46: was -795.3333
$ws.Cells.Item(46, 14).Value = -4873   # N46: was -6042.3335
$ws.Cells.Item(88, 8).Value = 0   # H88: was 18900
$ws.Cells.Item(88, 10).Value = 0   # J88: was 18900
$ws.Cells.Item(88, 12).Value = 0   # L88: was 18900
$ws.Cells.Item(88, 14).Value = ""   # N88: was -19756, clear cell
$ws.Cells.Item(91, 8).Value = 0   # H91: was 18900
$ws.Cells.Item(91, 10).Value = 0   # J91: was 18900
$ws.Cells.Item(91, 12).Value = 0   # L91: was 18900
$ws.Cells.Item(91, 14).Value = ""   # N91: was -21864, clear cell
$ws.Cells.Item(126, 8).Value = 21015.215   # H126: was 20425.068
$ws.Cells.Item(126, 9).Value = 44874.2   # I126: was 38045.082
$ws.Cells.Item(126, 10).Value = 7760.222   # J126: was 7987.4116
$ws.Cells.Item(126, 11).Value = 134622.6   # K126: was 114135.246
$ws.Cells.Item(126, 12).Value = 23280.666   # L126: was 23962.2348
$ws.Cells.Item(126, 13).Value = -132152.6   # M126: was -111665.246
$ws.Cells.Item(126, 14).Value = -28220.666   # N126: was -28902.2348
$ws.Cells.Item(132, 8).Value = 787865.2   # H132: was 576135.5600000001
$ws.Cells.Item(132, 9).Value = 1355949.9   # I132: was 1065709.9
$ws.Cells.Item(132, 10).Value = 6748.625   # J132: was 4965.5835
$ws.Cells.Item(132, 11).Value = 4067849.7   # K132: was 3197129.7
$ws.Cells.Item(132, 12).Value = 20245.875   # L132: was 14896.7505
$ws.Cells.Item(132, 13).Value = -4065319.7   # M132: was -3194599.7

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(30, 8).Value = 49000   # H30: was 25500
$ws.Cells.Item(30, 9).Value = 0   # I30: was 1000
$ws.Cells.Item(30, 10).Value = 49000   # J30: was 50000
$ws.Cells.Item(30, 11).Value = 0   # K30: was 1000
$ws.Cells.Item(30, 12).Value = 49000   # L30: was 50000
$ws.Cells.Item(30, 13).Value = ""   # M30: was -893, clear cell
$ws.Cells.Item(30, 14).Value = -49214   # N30: was -50214
$ws.Cells.Item(107, 8).Value = 11351   # H107: was 11752.893
$ws.Cells.Item(107, 9).Value = 1089.7084   # I107: was 1132.826
$ws.Cells.Item(107, 11).Value = 3269.1252   # K107: was 3398.478
$ws.Cells.Item(107, 13).Value = -1349.1252   # M107: was -1478.478
$ws.Cells.Item(122, 8).Value = 6483.8076   # H122: was 7216.522
$ws.Cells.Item(122, 9).Value = 2534.2856   # I122: was 2831.6667
$ws.Cells.Item(122, 10).Value = 11091.583   # J122: was 12000
$ws.Cells.Item(122, 11).Value = 7602.8568   # K122: was 8495.000100000001
$ws.Cells.Item(122, 12).Value = 33274.749   # L122: was 36000
$ws.Cells.Item(122, 13).Value = -5152.8568   # M122: was -6045.000100000001
$ws.Cells.Item(122, 14).Value = -38174.749   # N122: was -40900
$ws.Cells.Item(132, 8).Value = 8801.371999999999   # H132: was 8779.513000000001
$ws.Cells.Item(132, 9).Value = 9824.048000000001   # I132: was 9796.984
$ws.Cells.Item(132, 11).Value = 29472.144   # K132: was 29390.952
$ws.Cells.Item(132, 13).Value = -26942.144   # M132: was -26860.952
